# Update "paises" COVID data workbook (countries & Spain provincias update)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------------------
# Header timestamp cell
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 18:33"

# ---------------------------------------------------------------------------
# Row 4 - Estados Unidos: updated totals
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = 1243483
$ws.Range("C4").Value = 5850
$ws.Range("D4").Value = 202328
$ws.Range("E4").Value = 968326
$ws.Range("G4").Value = 558
$ws.Range("H4").Value = 72829

# ---------------------------------------------------------------------------
# Row 7 - Reino Unido: updated totals
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = 201101
$ws.Range("C7").Value = 6111
$ws.Range("E7").Value = 170681
$ws.Range("G7").Value = 649
$ws.Range("H7").Value = 30076

# ---------------------------------------------------------------------------
# Row 11 - Turquia: updated totals
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = 131744
$ws.Range("C11").Value = 2253
$ws.Range("D11").Value = 78202
$ws.Range("E11").Value = 49958
$ws.Range("F11").Value = 1278
$ws.Range("G11").Value = 64
$ws.Range("H11").Value = 3584

# ---------------------------------------------------------------------------
# Row 26/27 - Chile moves above Pakistan in the ranking, Chile gets fresh
# figures while Pakistan keeps its previous totals.
# ---------------------------------------------------------------------------
$ws.Range("A26").Value = "Chile"
$ws.Range("B26").Value = 23048
$ws.Range("C26").Value = 1032
$ws.Range("D26").Value = 11189
$ws.Range("E26").Value = 11578
$ws.Range("F26").Value = 470
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = 281

$ws.Range("A27").Value = "Pakistan"
$ws.Range("B27").Value = 22550
$ws.Range("C27").Value = 501
$ws.Range("D27").Value = 6217
$ws.Range("E27").Value = 15807
$ws.Range("F27").Value = 111
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = 526

# ---------------------------------------------------------------------------
# Row 47 - Noruega: updated totals
# ---------------------------------------------------------------------------
$ws.Range("B47").Value = 7976
$ws.Range("C47").Value = 21
$ws.Range("E47").Value = 7728
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 216

# ---------------------------------------------------------------------------
# Row 59 - Moldavia: updated totals
# ---------------------------------------------------------------------------
$ws.Range("B59").Value = 4476
$ws.Range("C59").Value = 113
$ws.Range("E59").Value = 2789
$ws.Range("G59").Value = 7
$ws.Range("H59").Value = 143

# ---------------------------------------------------------------------------
# Row 61/62 - Luxemburgo moves above Barein, Luxemburgo gets fresh figures
# while Barein keeps its previous totals.
# ---------------------------------------------------------------------------
$ws.Range("A61").Value = "Luxemburgo"
$ws.Range("B61").Value = 3851
$ws.Range("C61").Value = 11
$ws.Range("D61").Value = 3452
$ws.Range("E61").Value = 301
$ws.Range("F61").Value = 20
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 98

$ws.Range("A62").Value = "Barein"
$ws.Range("B62").Value = 3842
$ws.Range("C62").Value = 122
$ws.Range("D62").Value = 1860
$ws.Range("E62").Value = 1974
$ws.Range("F62").Value = 4
$ws.Range("H62").Value = 8

# ---------------------------------------------------------------------------
# Row 71 - Irak: updated totals
# ---------------------------------------------------------------------------
$ws.Range("B71").Value = 2480
$ws.Range("C71").Value = 49
$ws.Range("D71").Value = 1602
$ws.Range("E71").Value = 776

# ---------------------------------------------------------------------------
# Row 144 - Birmania: updated totals
# ---------------------------------------------------------------------------
$ws.Range("D144").Value = 50
$ws.Range("E144").Value = 105
